$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115..181 down to 116..182
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new record's data
$ws.Cells.Item(115, 1).Value2 = 3
$ws.Cells.Item(115, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(115, 3).Value2 = "Coquimbo"
$ws.Cells.Item(115, 4).Value2 = 44438
$ws.Cells.Item(115, 5).Value2 = 5
$ws.Cells.Item(115, 6).Value2 = 100112040
$ws.Cells.Item(115, 7).Value2 = "Cilantro"
$ws.Cells.Item(115, 8).Value2 = "Sin especificar"
$ws.Cells.Item(115, 9).Value2 = "Primera"
$ws.Cells.Item(115, 10).Value2 = 300
$ws.Cells.Item(115, 11).Value2 = 3300
$ws.Cells.Item(115, 12).Value2 = 3500
$ws.Cells.Item(115, 13).Value2 = 3420
$ws.Cells.Item(115, 14).Value2 = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(115, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(115, 16).Value2 = 1140
$ws.Cells.Item(115, 17).Value2 = 3
$ws.Cells.Item(115, 18).Value2 = "Hortaliza"
